$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly data point lands at row 176; every existing data row (old 176..291)
# shifts down by one (to 177..292), which is exactly what a native row insert does.
$ws.Rows.Item(176).Insert()

# Populate the freshly-inserted row with this week's record.
$ws.Range("A176").Value = 8
$ws.Range("B176").Value = "Terminal La Palmera de La Serena"
$ws.Range("C176").Value = "Coquimbo"
$ws.Range("D176").Value = 44762
$ws.Range("E176").Value = 4
$ws.Range("F176").Value = 100112012
$ws.Range("G176").Value = "Espinaca"
$ws.Range("H176").Value = "Sin especificar"
$ws.Range("I176").Value = "Primera"
$ws.Range("J176").Value = 2400
$ws.Range("K176").Value = 500
$ws.Range("L176").Value = 600
$ws.Range("M176").Value = 550
$ws.Range("N176").Value = "$/atado 300 a 500 gramos"
$ws.Range("O176").Value = "Provincia del Elquí"
$ws.Range("P176").Value = 1100
$ws.Range("Q176").Value = 0.5
$ws.Range("R176").Value = "Hortaliza"
